# ANALISIS DE LA IMAGEN RANA
# Adds the "RUIDO" analysis mini-table (K3:M8, with a merged comment-style
# note box) and a second "PSD" table (rows 11-15), fills in the missing
# D8/E8 pair ("INFINITO" sample), removes the now-unused explicit style
# from C5:C8/F5:G5, drops the two empty trailing helper rows 9-10, and
# attaches the threaded review comment on K4.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- complete the last row of the first table -----------------------------
$ws.Range("D8").Value = 32
$ws.Range("E8").Value = "INFINITO"

# --- drop the stray formatting left on the numeric/formula cells ----------
$ws.Range("C5:C8").ClearFormats()
$ws.Range("F5:G5").ClearFormats()

# --- remove the two leftover empty rows ------------------------------------
$ws.Range("C9").Clear()
$ws.Range("C10").Clear()

# --- new "RUIDO" table (K3:M8) --------------------------------------------
$ws.Range("K3").Value = "RUIDO "

$ws.Range("K4").Value = "PSD AMPLITUD "
$ws.Range("L4").Value = "SPS"
$ws.Range("M4").Value = "TIEMPO "

$ws.Range("K5").Value = 0.0000313
$ws.Range("L5").Value = 1
$ws.Range("M5").Value = "AUMENTA LA FRECUNCIA DE ALETORIEDAD "

$ws.Range("K6").Value = 0.00000775
$ws.Range("L6").Value = 4

$ws.Range("K7").Value = 0.0000039
$ws.Range("L7").Value = 8

$ws.Range("K8").Value = 0.00000196
$ws.Range("L8").Value = 16

$ws.Range("K5:K8").NumberFormat = "0.00E+00"

$ws.Range("M5:O8").Merge()
$ws.Range("M5:O8").HorizontalAlignment = -4108
$ws.Range("M5:O8").VerticalAlignment = -4108
$ws.Range("M5:O8").WrapText = $true

$ws.Columns("K:K").ColumnWidth = 16.109375

# --- threaded review comment on K4 -----------------------------------------
$ws.Range("K4").AddCommentThreaded("A medida que tomamos mas muestras por segundo se puede apreciar mas el comportamiento aleatorio del ruido con mas frecuencia y adicional disminuye la potencia y se hace mas continua. ")

# --- second table: PSD vs Sps ----------------------------------------------
$ws.Range("C11").Value = "PSD "

$ws.Range("B12").Value = 1
$ws.Range("C12").Value = 0.0008

$ws.Range("B13").Value = 4
$ws.Range("C13").Value = 0.0006

$ws.Range("B14").Value = 8
$ws.Range("C14").Value = 0.00036

$ws.Range("B15").Value = 16
$ws.Range("C15").Value = 0.00022

# --- selection matches the saved state --------------------------------------
$ws.Range("H21").Select()
